$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns before D, shifting existing D:K data to F:M
$ws.Range("D:E").Insert()

# Step 2: fix up styles for the two new columns (Insert leaves them at the default
# left-neighbour style; copy the number-format/font from column F of each row so the
# new D/E cells match the date-header style (2) or numeric style (3) used across the row)
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("D8:E8").PasteSpecial(-4122)
$ws.Range("F9").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)
$ws.Range("F10").Copy()
$ws.Range("D10:E10").PasteSpecial(-4122)
$ws.Range("F11").Copy()
$ws.Range("D11:E11").PasteSpecial(-4122)
$ws.Range("F12").Copy()
$ws.Range("D12:E12").PasteSpecial(-4122)
$ws.Range("F13").Copy()
$ws.Range("D13:E13").PasteSpecial(-4122)
$ws.Range("F14").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)
$ws.Range("F15").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)
$ws.Range("F16").Copy()
$ws.Range("D16:E16").PasteSpecial(-4122)
$ws.Range("F17").Copy()
$ws.Range("D17:E17").PasteSpecial(-4122)
$ws.Range("F18").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)
$ws.Range("F19").Copy()
$ws.Range("D19:E19").PasteSpecial(-4122)
$ws.Range("F20").Copy()
$ws.Range("D20:E20").PasteSpecial(-4122)
$ws.Range("F21").Copy()
$ws.Range("D21:E21").PasteSpecial(-4122)
$ws.Range("F22").Copy()
$ws.Range("D22:E22").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("D23:E23").PasteSpecial(-4122)
$ws.Range("F24").Copy()
$ws.Range("D24:E24").PasteSpecial(-4122)
$ws.Range("F25").Copy()
$ws.Range("D25:E25").PasteSpecial(-4122)
$ws.Range("F26").Copy()
$ws.Range("D26:E26").PasteSpecial(-4122)
$ws.Range("F27").Copy()
$ws.Range("D27:E27").PasteSpecial(-4122)
$ws.Range("F28").Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
$ws.Range("F29").Copy()
$ws.Range("D29:E29").PasteSpecial(-4122)
$ws.Range("F30").Copy()
$ws.Range("D30:E30").PasteSpecial(-4122)
$ws.Range("F31").Copy()
$ws.Range("D31:E31").PasteSpecial(-4122)
$ws.Range("F32").Copy()
$ws.Range("D32:E32").PasteSpecial(-4122)
$ws.Range("F33").Copy()
$ws.Range("D33:E33").PasteSpecial(-4122)
$ws.Range("F34").Copy()
$ws.Range("D34:E34").PasteSpecial(-4122)
$ws.Range("F35").Copy()
$ws.Range("D35:E35").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F39").Copy()
$ws.Range("D39:E39").PasteSpecial(-4122)
$ws.Range("F40").Copy()
$ws.Range("D40:E40").PasteSpecial(-4122)
$ws.Range("F41").Copy()
$ws.Range("D41:E41").PasteSpecial(-4122)
$ws.Range("F42").Copy()
$ws.Range("D42:E42").PasteSpecial(-4122)
$ws.Range("F43").Copy()
$ws.Range("D43:E43").PasteSpecial(-4122)
$ws.Range("F44").Copy()
$ws.Range("D44:E44").PasteSpecial(-4122)
$ws.Range("F45").Copy()
$ws.Range("D45:E45").PasteSpecial(-4122)
$ws.Range("F46").Copy()
$ws.Range("D46:E46").PasteSpecial(-4122)
$ws.Range("F47").Copy()
$ws.Range("D47:E47").PasteSpecial(-4122)
$ws.Range("F48").Copy()
$ws.Range("D48:E48").PasteSpecial(-4122)
$ws.Range("F49").Copy()
$ws.Range("D49:E49").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("D50:E50").PasteSpecial(-4122)
$ws.Range("F51").Copy()
$ws.Range("D51:E51").PasteSpecial(-4122)
$ws.Range("F52").Copy()
$ws.Range("D52:E52").PasteSpecial(-4122)
$ws.Range("F53").Copy()
$ws.Range("D53:E53").PasteSpecial(-4122)
$ws.Range("F54").Copy()
$ws.Range("D54:E54").PasteSpecial(-4122)
$ws.Range("F55").Copy()
$ws.Range("D55:E55").PasteSpecial(-4122)
$ws.Range("F56").Copy()
$ws.Range("D56:E56").PasteSpecial(-4122)
$ws.Range("F57").Copy()
$ws.Range("D57:E57").PasteSpecial(-4122)
$ws.Range("F58").Copy()
$ws.Range("D58:E58").PasteSpecial(-4122)
$ws.Range("F59").Copy()
$ws.Range("D59:E59").PasteSpecial(-4122)
$ws.Range("F60").Copy()
$ws.Range("D60:E60").PasteSpecial(-4122)
$ws.Range("F61").Copy()
$ws.Range("D61:E61").PasteSpecial(-4122)
$ws.Range("F62").Copy()
$ws.Range("D62:E62").PasteSpecial(-4122)
$ws.Range("F63").Copy()
$ws.Range("D63:E63").PasteSpecial(-4122)
$ws.Range("F64").Copy()
$ws.Range("D64:E64").PasteSpecial(-4122)
$ws.Range("F65").Copy()
$ws.Range("D65:E65").PasteSpecial(-4122)
$ws.Range("F66").Copy()
$ws.Range("D66:E66").PasteSpecial(-4122)
$ws.Range("F67").Copy()
$ws.Range("D67:E67").PasteSpecial(-4122)
$ws.Range("F68").Copy()
$ws.Range("D68:E68").PasteSpecial(-4122)
$ws.Range("F69").Copy()
$ws.Range("D69:E69").PasteSpecial(-4122)
$ws.Range("F70").Copy()
$ws.Range("D70:E70").PasteSpecial(-4122)
$ws.Range("F71").Copy()
$ws.Range("D71:E71").PasteSpecial(-4122)
$ws.Range("F72").Copy()
$ws.Range("D72:E72").PasteSpecial(-4122)
$ws.Range("F73").Copy()
$ws.Range("D73:E73").PasteSpecial(-4122)
$ws.Range("F74").Copy()
$ws.Range("D74:E74").PasteSpecial(-4122)
$ws.Range("F75").Copy()
$ws.Range("D75:E75").PasteSpecial(-4122)
$ws.Range("F76").Copy()
$ws.Range("D76:E76").PasteSpecial(-4122)
$ws.Range("F77").Copy()
$ws.Range("D77:E77").PasteSpecial(-4122)
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Range("F81").Copy()
$ws.Range("D81:E81").PasteSpecial(-4122)
$ws.Range("F82").Copy()
$ws.Range("D82:E82").PasteSpecial(-4122)
$ws.Range("F83").Copy()
$ws.Range("D83:E83").PasteSpecial(-4122)
$ws.Range("F84").Copy()
$ws.Range("D84:E84").PasteSpecial(-4122)
$ws.Range("F85").Copy()
$ws.Range("D85:E85").PasteSpecial(-4122)
$ws.Range("F86").Copy()
$ws.Range("D86:E86").PasteSpecial(-4122)
$ws.Range("F87").Copy()
$ws.Range("D87:E87").PasteSpecial(-4122)
$ws.Range("F88").Copy()
$ws.Range("D88:E88").PasteSpecial(-4122)
$ws.Range("F89").Copy()
$ws.Range("D89:E89").PasteSpecial(-4122)
$ws.Range("F90").Copy()
$ws.Range("D90:E90").PasteSpecial(-4122)
$ws.Range("F91").Copy()
$ws.Range("D91:E91").PasteSpecial(-4122)
$ws.Range("F92").Copy()
$ws.Range("D92:E92").PasteSpecial(-4122)
$ws.Range("F93").Copy()
$ws.Range("D93:E93").PasteSpecial(-4122)
$ws.Range("F94").Copy()
$ws.Range("D94:E94").PasteSpecial(-4122)
$ws.Range("F95").Copy()
$ws.Range("D95:E95").PasteSpecial(-4122)
$ws.Range("F96").Copy()
$ws.Range("D96:E96").PasteSpecial(-4122)
$ws.Range("F97").Copy()
$ws.Range("D97:E97").PasteSpecial(-4122)
$ws.Range("F98").Copy()
$ws.Range("D98:E98").PasteSpecial(-4122)
$ws.Range("F99").Copy()
$ws.Range("D99:E99").PasteSpecial(-4122)
$ws.Range("F100").Copy()
$ws.Range("D100:E100").PasteSpecial(-4122)
$ws.Range("F101").Copy()
$ws.Range("D101:E101").PasteSpecial(-4122)
$ws.Range("F102").Copy()
$ws.Range("D102:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: write the refreshed financial figures (columns D:M) for every data row
$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(7,5).Value = 43373
$ws.Cells.Item(7,6).Value = 43281
$ws.Cells.Item(7,7).Value = 43190
$ws.Cells.Item(7,8).Value = 43100
$ws.Cells.Item(7,9).Value = 43008
$ws.Cells.Item(7,10).Value = 42916
$ws.Cells.Item(7,11).Value = 42825
$ws.Cells.Item(7,12).Value = 42735
$ws.Cells.Item(7,13).Value = 42643
$ws.Cells.Item(8,4).Value = 406800
$ws.Cells.Item(8,5).Value = 282200
$ws.Cells.Item(8,6).Value = 218500
$ws.Cells.Item(8,7).Value = 181500
$ws.Cells.Item(8,8).Value = 168700
$ws.Cells.Item(8,9).Value = 145100
$ws.Cells.Item(8,10).Value = 129300
$ws.Cells.Item(8,11).Value = 138300
$ws.Cells.Item(8,12).Value = 145100
$ws.Cells.Item(8,13).Value = 143300
$ws.Cells.Item(9,4).Value = "NA"
$ws.Cells.Item(9,5).Value = "NA"
$ws.Cells.Item(9,6).Value = "NA"
$ws.Cells.Item(9,7).Value = "NA"
$ws.Cells.Item(9,8).Value = "NA"
$ws.Cells.Item(9,9).Value = "NA"
$ws.Cells.Item(9,10).Value = "NA"
$ws.Cells.Item(9,11).Value = "NA"
$ws.Cells.Item(9,12).Value = "NA"
$ws.Cells.Item(9,13).Value = "NA"
$ws.Cells.Item(10,4).Value = "NA"
$ws.Cells.Item(10,5).Value = "NA"
$ws.Cells.Item(10,6).Value = "NA"
$ws.Cells.Item(10,7).Value = "NA"
$ws.Cells.Item(10,8).Value = "NA"
$ws.Cells.Item(10,9).Value = "NA"
$ws.Cells.Item(10,10).Value = "NA"
$ws.Cells.Item(10,11).Value = "NA"
$ws.Cells.Item(10,12).Value = "NA"
$ws.Cells.Item(10,13).Value = "NA"
$ws.Cells.Item(12,4).Value = "NA"
$ws.Cells.Item(12,5).Value = "NA"
$ws.Cells.Item(12,6).Value = "NA"
$ws.Cells.Item(12,7).Value = "NA"
$ws.Cells.Item(12,8).Value = "NA"
$ws.Cells.Item(12,9).Value = "NA"
$ws.Cells.Item(12,10).Value = "NA"
$ws.Cells.Item(12,11).Value = "NA"
$ws.Cells.Item(12,12).Value = "NA"
$ws.Cells.Item(12,13).Value = "NA"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 0
$ws.Cells.Item(13,8).Value = 0
$ws.Cells.Item(13,9).Value = 0
$ws.Cells.Item(13,10).Value = 0
$ws.Cells.Item(13,11).Value = 0
$ws.Cells.Item(13,12).Value = 0
$ws.Cells.Item(13,13).Value = 0
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = 0
$ws.Cells.Item(14,9).Value = 0
$ws.Cells.Item(14,10).Value = 0
$ws.Cells.Item(14,11).Value = 0
$ws.Cells.Item(14,12).Value = 0
$ws.Cells.Item(14,13).Value = 0
$ws.Cells.Item(15,4).Value = -6400
$ws.Cells.Item(15,5).Value = -5000
$ws.Cells.Item(15,6).Value = -4800
$ws.Cells.Item(15,7).Value = -4600
$ws.Cells.Item(15,8).Value = -4100
$ws.Cells.Item(15,9).Value = -3800
$ws.Cells.Item(15,10).Value = -3300
$ws.Cells.Item(15,11).Value = -3000
$ws.Cells.Item(15,12).Value = -2500
$ws.Cells.Item(15,13).Value = -2200
$ws.Cells.Item(17,4).Value = 237800
$ws.Cells.Item(17,5).Value = 155300
$ws.Cells.Item(17,6).Value = 98900
$ws.Cells.Item(17,7).Value = 76500
$ws.Cells.Item(17,8).Value = 72400
$ws.Cells.Item(17,9).Value = 65500
$ws.Cells.Item(17,10).Value = 58700
$ws.Cells.Item(17,11).Value = 66500
$ws.Cells.Item(17,12).Value = 68600
$ws.Cells.Item(17,13).Value = 74900
$ws.Cells.Item(18,4).Value = 169000
$ws.Cells.Item(18,5).Value = 126900
$ws.Cells.Item(18,6).Value = 119600
$ws.Cells.Item(18,7).Value = 105000
$ws.Cells.Item(18,8).Value = 96300
$ws.Cells.Item(18,9).Value = 79600
$ws.Cells.Item(18,10).Value = 70600
$ws.Cells.Item(18,11).Value = 71800
$ws.Cells.Item(18,12).Value = 76500
$ws.Cells.Item(18,13).Value = 68400
$ws.Cells.Item(20,4).Value = -70200
$ws.Cells.Item(20,5).Value = -31000
$ws.Cells.Item(20,6).Value = -50300
$ws.Cells.Item(20,7).Value = -53900
$ws.Cells.Item(20,8).Value = -46500
$ws.Cells.Item(20,9).Value = -37100
$ws.Cells.Item(20,10).Value = -66000
$ws.Cells.Item(20,11).Value = -21500
$ws.Cells.Item(20,12).Value = -50400
$ws.Cells.Item(20,13).Value = -28400
$ws.Cells.Item(21,4).Value = 104500
$ws.Cells.Item(21,5).Value = 100900
$ws.Cells.Item(21,6).Value = 74100
$ws.Cells.Item(21,7).Value = 55700
$ws.Cells.Item(21,8).Value = 50000
$ws.Cells.Item(21,9).Value = 42900
$ws.Cells.Item(21,10).Value = 5100
$ws.Cells.Item(21,11).Value = 52600
$ws.Cells.Item(21,12).Value = "NA"
$ws.Cells.Item(21,13).Value = 40200
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 0
$ws.Cells.Item(22,8).Value = 0
$ws.Cells.Item(22,9).Value = 0
$ws.Cells.Item(22,10).Value = 0
$ws.Cells.Item(22,11).Value = 0
$ws.Cells.Item(22,12).Value = 0
$ws.Cells.Item(22,13).Value = 0
$ws.Cells.Item(23,4).Value = 98700
$ws.Cells.Item(23,5).Value = 95900
$ws.Cells.Item(23,6).Value = 69300
$ws.Cells.Item(23,7).Value = 51100
$ws.Cells.Item(23,8).Value = 49800
$ws.Cells.Item(23,9).Value = 42500
$ws.Cells.Item(23,10).Value = 4500
$ws.Cells.Item(23,11).Value = 50300
$ws.Cells.Item(23,12).Value = 26200
$ws.Cells.Item(23,13).Value = 40000
$ws.Cells.Item(24,4).Value = 31300
$ws.Cells.Item(24,5).Value = 26000
$ws.Cells.Item(24,6).Value = 19900
$ws.Cells.Item(24,7).Value = 15200
$ws.Cells.Item(24,8).Value = 17200
$ws.Cells.Item(24,9).Value = 16400
$ws.Cells.Item(24,10).Value = 11000
$ws.Cells.Item(24,11).Value = 8300
$ws.Cells.Item(24,12).Value = 10900
$ws.Cells.Item(24,13).Value = 15600
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = 0
$ws.Cells.Item(25,9).Value = 0
$ws.Cells.Item(25,10).Value = 0
$ws.Cells.Item(25,11).Value = 0
$ws.Cells.Item(25,12).Value = 0
$ws.Cells.Item(25,13).Value = 0
$ws.Cells.Item(26,4).Value = 67400
$ws.Cells.Item(26,5).Value = 69900
$ws.Cells.Item(26,6).Value = 49400
$ws.Cells.Item(26,7).Value = 35900
$ws.Cells.Item(26,8).Value = 32600
$ws.Cells.Item(26,9).Value = 26000
$ws.Cells.Item(26,10).Value = -6500
$ws.Cells.Item(26,11).Value = 42000
$ws.Cells.Item(26,12).Value = 15200
$ws.Cells.Item(26,13).Value = 24400
$ws.Cells.Item(27,4).Value = 67400
$ws.Cells.Item(27,5).Value = 69000
$ws.Cells.Item(27,6).Value = 48700
$ws.Cells.Item(27,7).Value = 35400
$ws.Cells.Item(27,8).Value = 32600
$ws.Cells.Item(27,9).Value = 26000
$ws.Cells.Item(27,10).Value = -6500
$ws.Cells.Item(27,11).Value = 42000
$ws.Cells.Item(27,12).Value = 15200
$ws.Cells.Item(27,13).Value = 24400
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = 0
$ws.Cells.Item(28,7).Value = 0
$ws.Cells.Item(28,8).Value = 0
$ws.Cells.Item(28,9).Value = 0
$ws.Cells.Item(28,10).Value = 0
$ws.Cells.Item(28,11).Value = 0
$ws.Cells.Item(28,12).Value = 0
$ws.Cells.Item(28,13).Value = 0
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,6).Value = 0
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = 0
$ws.Cells.Item(29,9).Value = 0
$ws.Cells.Item(29,10).Value = 0
$ws.Cells.Item(29,11).Value = 0
$ws.Cells.Item(29,12).Value = 0
$ws.Cells.Item(29,13).Value = 0
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 0
$ws.Cells.Item(30,7).Value = 0
$ws.Cells.Item(30,8).Value = 0
$ws.Cells.Item(30,9).Value = 0
$ws.Cells.Item(30,10).Value = 0
$ws.Cells.Item(30,11).Value = 0
$ws.Cells.Item(30,12).Value = 0
$ws.Cells.Item(30,13).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(31,6).Value = 0
$ws.Cells.Item(31,7).Value = 0
$ws.Cells.Item(31,8).Value = 0
$ws.Cells.Item(31,9).Value = 0
$ws.Cells.Item(31,10).Value = 0
$ws.Cells.Item(31,11).Value = 0
$ws.Cells.Item(31,12).Value = 0
$ws.Cells.Item(31,13).Value = 0
$ws.Cells.Item(32,4).Value = 70200
$ws.Cells.Item(32,5).Value = 31000
$ws.Cells.Item(32,6).Value = 50300
$ws.Cells.Item(32,7).Value = 53900
$ws.Cells.Item(32,8).Value = 46500
$ws.Cells.Item(32,9).Value = 37100
$ws.Cells.Item(32,10).Value = 66000
$ws.Cells.Item(32,11).Value = 21500
$ws.Cells.Item(32,12).Value = 50400
$ws.Cells.Item(32,13).Value = 28400
$ws.Cells.Item(33,4).Value = 67400
$ws.Cells.Item(33,5).Value = 69000
$ws.Cells.Item(33,6).Value = 48700
$ws.Cells.Item(33,7).Value = 35400
$ws.Cells.Item(33,8).Value = 32600
$ws.Cells.Item(33,9).Value = 26000
$ws.Cells.Item(33,10).Value = -6500
$ws.Cells.Item(33,11).Value = 42000
$ws.Cells.Item(33,12).Value = 15200
$ws.Cells.Item(33,13).Value = 24400
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(34,6).Value = 0
$ws.Cells.Item(34,7).Value = 0
$ws.Cells.Item(34,8).Value = 0
$ws.Cells.Item(34,9).Value = 0
$ws.Cells.Item(34,10).Value = 0
$ws.Cells.Item(34,11).Value = 0
$ws.Cells.Item(34,12).Value = 0
$ws.Cells.Item(34,13).Value = 0
$ws.Cells.Item(35,4).Value = 67400
$ws.Cells.Item(35,5).Value = 69000
$ws.Cells.Item(35,6).Value = 48700
$ws.Cells.Item(35,7).Value = 35400
$ws.Cells.Item(35,8).Value = 32600
$ws.Cells.Item(35,9).Value = 26000
$ws.Cells.Item(35,10).Value = -6500
$ws.Cells.Item(35,11).Value = 42000
$ws.Cells.Item(35,12).Value = 15200
$ws.Cells.Item(35,13).Value = 24400
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(38,5).Value = 43373
$ws.Cells.Item(38,6).Value = 43281
$ws.Cells.Item(38,7).Value = 43190
$ws.Cells.Item(38,8).Value = 43100
$ws.Cells.Item(38,9).Value = 43008
$ws.Cells.Item(38,10).Value = 42916
$ws.Cells.Item(38,11).Value = 42825
$ws.Cells.Item(38,12).Value = 42735
$ws.Cells.Item(38,13).Value = 42643
$ws.Cells.Item(41,4).Value = 2089300
$ws.Cells.Item(41,5).Value = 1865700
$ws.Cells.Item(41,6).Value = 1153100
$ws.Cells.Item(41,7).Value = 743300
$ws.Cells.Item(41,8).Value = 828300
$ws.Cells.Item(41,9).Value = 535600
$ws.Cells.Item(41,10).Value = 839500
$ws.Cells.Item(41,11).Value = 923100
$ws.Cells.Item(41,12).Value = 1189700
$ws.Cells.Item(41,13).Value = 763300
$ws.Cells.Item(42,4).Value = 359400
$ws.Cells.Item(42,5).Value = 162800
$ws.Cells.Item(42,6).Value = 163400
$ws.Cells.Item(42,7).Value = 133700
$ws.Cells.Item(42,8).Value = 387200
$ws.Cells.Item(42,9).Value = 441500
$ws.Cells.Item(42,10).Value = 171400
$ws.Cells.Item(42,11).Value = 302800
$ws.Cells.Item(42,12).Value = 185500
$ws.Cells.Item(42,13).Value = 116900
$ws.Cells.Item(43,4).Value = 0
$ws.Cells.Item(43,5).Value = 0
$ws.Cells.Item(43,6).Value = 0
$ws.Cells.Item(43,7).Value = 0
$ws.Cells.Item(43,8).Value = 0
$ws.Cells.Item(43,9).Value = 0
$ws.Cells.Item(43,10).Value = 0
$ws.Cells.Item(43,11).Value = 0
$ws.Cells.Item(43,12).Value = 0
$ws.Cells.Item(43,13).Value = 0
$ws.Cells.Item(44,4).Value = 0
$ws.Cells.Item(44,5).Value = 0
$ws.Cells.Item(44,6).Value = 0
$ws.Cells.Item(44,7).Value = 0
$ws.Cells.Item(44,8).Value = 0
$ws.Cells.Item(44,9).Value = 0
$ws.Cells.Item(44,10).Value = 0
$ws.Cells.Item(44,11).Value = 0
$ws.Cells.Item(44,12).Value = 0
$ws.Cells.Item(44,13).Value = 0
$ws.Cells.Item(45,4).Value = 0
$ws.Cells.Item(45,5).Value = 0
$ws.Cells.Item(45,6).Value = 0
$ws.Cells.Item(45,7).Value = 0
$ws.Cells.Item(45,8).Value = 0
$ws.Cells.Item(45,9).Value = 0
$ws.Cells.Item(45,10).Value = 0
$ws.Cells.Item(45,11).Value = 0
$ws.Cells.Item(45,12).Value = 0
$ws.Cells.Item(45,13).Value = 0
$ws.Cells.Item(46,4).Value = 0
$ws.Cells.Item(46,5).Value = 0
$ws.Cells.Item(46,6).Value = 0
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(46,8).Value = 0
$ws.Cells.Item(46,9).Value = 0
$ws.Cells.Item(46,10).Value = 0
$ws.Cells.Item(46,11).Value = 0
$ws.Cells.Item(46,12).Value = 0
$ws.Cells.Item(46,13).Value = 0
$ws.Cells.Item(47,4).Value = 40200
$ws.Cells.Item(47,5).Value = 40900
$ws.Cells.Item(47,6).Value = 25000
$ws.Cells.Item(47,7).Value = 21100
$ws.Cells.Item(47,8).Value = 15100
$ws.Cells.Item(47,9).Value = 13800
$ws.Cells.Item(47,10).Value = 11700
$ws.Cells.Item(47,11).Value = 14000
$ws.Cells.Item(47,12).Value = 13300
$ws.Cells.Item(47,13).Value = 12900
$ws.Cells.Item(48,4).Value = 225200
$ws.Cells.Item(48,5).Value = 210400
$ws.Cells.Item(48,6).Value = 209600
$ws.Cells.Item(48,7).Value = 207000
$ws.Cells.Item(48,8).Value = 97400
$ws.Cells.Item(48,9).Value = 92400
$ws.Cells.Item(48,10).Value = 88500
$ws.Cells.Item(48,11).Value = 87000
$ws.Cells.Item(48,12).Value = 83700
$ws.Cells.Item(48,13).Value = 31800
$ws.Cells.Item(49,4).Value = 11700
$ws.Cells.Item(49,5).Value = 13700
$ws.Cells.Item(49,6).Value = 12400
$ws.Cells.Item(49,7).Value = 10300
$ws.Cells.Item(49,8).Value = 10000
$ws.Cells.Item(49,9).Value = 9400
$ws.Cells.Item(49,10).Value = 8100
$ws.Cells.Item(49,11).Value = 8700
$ws.Cells.Item(49,12).Value = 8300
$ws.Cells.Item(49,13).Value = 7500
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(50,6).Value = 0
$ws.Cells.Item(50,7).Value = 0
$ws.Cells.Item(50,8).Value = 0
$ws.Cells.Item(50,9).Value = 0
$ws.Cells.Item(50,10).Value = 0
$ws.Cells.Item(50,11).Value = 0
$ws.Cells.Item(50,12).Value = 0
$ws.Cells.Item(50,13).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(51,6).Value = 0
$ws.Cells.Item(51,7).Value = 0
$ws.Cells.Item(51,8).Value = 0
$ws.Cells.Item(51,9).Value = 0
$ws.Cells.Item(51,10).Value = 0
$ws.Cells.Item(51,11).Value = 0
$ws.Cells.Item(51,12).Value = 0
$ws.Cells.Item(51,13).Value = 0
$ws.Cells.Item(52,4).Value = 15800
$ws.Cells.Item(52,5).Value = 8100
$ws.Cells.Item(52,6).Value = 12200
$ws.Cells.Item(52,7).Value = 18300
$ws.Cells.Item(52,8).Value = 18200
$ws.Cells.Item(52,9).Value = 19700
$ws.Cells.Item(52,10).Value = 17200
$ws.Cells.Item(52,11).Value = 16200
$ws.Cells.Item(52,12).Value = "NA"
$ws.Cells.Item(52,13).Value = "NA"
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(53,6).Value = 0
$ws.Cells.Item(53,7).Value = 0
$ws.Cells.Item(53,8).Value = 0
$ws.Cells.Item(53,9).Value = 0
$ws.Cells.Item(53,10).Value = 0
$ws.Cells.Item(53,11).Value = 0
$ws.Cells.Item(53,12).Value = 0
$ws.Cells.Item(53,13).Value = 0
$ws.Cells.Item(54,4).Value = 8134800
$ws.Cells.Item(54,5).Value = 7661600
$ws.Cells.Item(54,6).Value = 6153200
$ws.Cells.Item(54,7).Value = 5144900
$ws.Cells.Item(54,8).Value = 5176200
$ws.Cells.Item(54,9).Value = 4589300
$ws.Cells.Item(54,10).Value = 3956000
$ws.Cells.Item(54,11).Value = 4826700
$ws.Cells.Item(54,12).Value = 3972900
$ws.Cells.Item(54,13).Value = 3427800
$ws.Cells.Item(57,4).Value = 673200
$ws.Cells.Item(57,5).Value = 535000
$ws.Cells.Item(57,6).Value = 445500
$ws.Cells.Item(57,7).Value = 418200
$ws.Cells.Item(57,8).Value = 382500
$ws.Cells.Item(57,9).Value = 392300
$ws.Cells.Item(57,10).Value = 276200
$ws.Cells.Item(57,11).Value = 453100
$ws.Cells.Item(57,12).Value = 266200
$ws.Cells.Item(57,13).Value = 278600
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 0
$ws.Cells.Item(58,6).Value = 0
$ws.Cells.Item(58,7).Value = 0
$ws.Cells.Item(58,8).Value = 0
$ws.Cells.Item(58,9).Value = 0
$ws.Cells.Item(58,10).Value = 0
$ws.Cells.Item(58,11).Value = 0
$ws.Cells.Item(58,12).Value = 0
$ws.Cells.Item(58,13).Value = 0
$ws.Cells.Item(59,4).Value = 164700
$ws.Cells.Item(59,5).Value = 135200
$ws.Cells.Item(59,6).Value = 90200
$ws.Cells.Item(59,7).Value = 69200
$ws.Cells.Item(59,8).Value = 75200
$ws.Cells.Item(59,9).Value = 59900
$ws.Cells.Item(59,10).Value = 65900
$ws.Cells.Item(59,11).Value = 56700
$ws.Cells.Item(59,12).Value = 67100
$ws.Cells.Item(59,13).Value = 62900
$ws.Cells.Item(60,4).Value = 0
$ws.Cells.Item(60,5).Value = 0
$ws.Cells.Item(60,6).Value = 0
$ws.Cells.Item(60,7).Value = 0
$ws.Cells.Item(60,8).Value = 0
$ws.Cells.Item(60,9).Value = 0
$ws.Cells.Item(60,10).Value = 0
$ws.Cells.Item(60,11).Value = 0
$ws.Cells.Item(60,12).Value = 0
$ws.Cells.Item(60,13).Value = 0
$ws.Cells.Item(61,4).Value = 300
$ws.Cells.Item(61,5).Value = 0
$ws.Cells.Item(61,6).Value = 21500
$ws.Cells.Item(61,7).Value = 13300
$ws.Cells.Item(61,8).Value = 108700
$ws.Cells.Item(61,9).Value = 61700
$ws.Cells.Item(61,10).Value = 61400
$ws.Cells.Item(61,11).Value = 69700
$ws.Cells.Item(61,12).Value = 83000
$ws.Cells.Item(61,13).Value = 77500
$ws.Cells.Item(62,4).Value = 84400
$ws.Cells.Item(62,5).Value = 80200
$ws.Cells.Item(62,6).Value = 85300
$ws.Cells.Item(62,7).Value = 86800
$ws.Cells.Item(62,8).Value = 64300
$ws.Cells.Item(62,9).Value = 64300
$ws.Cells.Item(62,10).Value = 60000
$ws.Cells.Item(62,11).Value = 35300
$ws.Cells.Item(62,12).Value = 36800
$ws.Cells.Item(62,13).Value = 31400
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(63,6).Value = 0
$ws.Cells.Item(63,7).Value = 0
$ws.Cells.Item(63,8).Value = 0
$ws.Cells.Item(63,9).Value = 0
$ws.Cells.Item(63,10).Value = 0
$ws.Cells.Item(63,11).Value = 0
$ws.Cells.Item(63,12).Value = 0
$ws.Cells.Item(63,13).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(64,6).Value = 0
$ws.Cells.Item(64,7).Value = 0
$ws.Cells.Item(64,8).Value = 0
$ws.Cells.Item(64,9).Value = 0
$ws.Cells.Item(64,10).Value = 0
$ws.Cells.Item(64,11).Value = 0
$ws.Cells.Item(64,12).Value = 0
$ws.Cells.Item(64,13).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(65,6).Value = 0
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(65,8).Value = 0
$ws.Cells.Item(65,9).Value = 0
$ws.Cells.Item(65,10).Value = 0
$ws.Cells.Item(65,11).Value = 0
$ws.Cells.Item(65,12).Value = 0
$ws.Cells.Item(65,13).Value = 0
$ws.Cells.Item(66,4).Value = 7250500
$ws.Cells.Item(66,5).Value = 6842600
$ws.Cells.Item(66,6).Value = 5399300
$ws.Cells.Item(66,7).Value = 4418900
$ws.Cells.Item(66,8).Value = 4578500
$ws.Cells.Item(66,9).Value = 4024100
$ws.Cells.Item(66,10).Value = 3568900
$ws.Cells.Item(66,11).Value = 4377500
$ws.Cells.Item(66,12).Value = 3542000
$ws.Cells.Item(66,13).Value = 3012100
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(68,6).Value = 0
$ws.Cells.Item(68,7).Value = 0
$ws.Cells.Item(68,8).Value = 0
$ws.Cells.Item(68,9).Value = 0
$ws.Cells.Item(68,10).Value = 0
$ws.Cells.Item(68,11).Value = 0
$ws.Cells.Item(68,12).Value = 0
$ws.Cells.Item(68,13).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(69,6).Value = 0
$ws.Cells.Item(69,7).Value = 0
$ws.Cells.Item(69,8).Value = 0
$ws.Cells.Item(69,9).Value = 0
$ws.Cells.Item(69,10).Value = 0
$ws.Cells.Item(69,11).Value = 0
$ws.Cells.Item(69,12).Value = 0
$ws.Cells.Item(69,13).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(70,6).Value = 0
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = 0
$ws.Cells.Item(70,9).Value = 0
$ws.Cells.Item(70,10).Value = 0
$ws.Cells.Item(70,11).Value = 0
$ws.Cells.Item(70,12).Value = 0
$ws.Cells.Item(70,13).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(71,6).Value = 0
$ws.Cells.Item(71,7).Value = 0
$ws.Cells.Item(71,8).Value = 0
$ws.Cells.Item(71,9).Value = 0
$ws.Cells.Item(71,10).Value = 0
$ws.Cells.Item(71,11).Value = 0
$ws.Cells.Item(71,12).Value = 0
$ws.Cells.Item(71,13).Value = 0
$ws.Cells.Item(72,4).Value = 708600
$ws.Cells.Item(72,5).Value = 643200
$ws.Cells.Item(72,6).Value = 578100
$ws.Cells.Item(72,7).Value = 550200
$ws.Cells.Item(72,8).Value = 422000
$ws.Cells.Item(72,9).Value = 389400
$ws.Cells.Item(72,10).Value = 363400
$ws.Cells.Item(72,11).Value = 422100
$ws.Cells.Item(72,12).Value = 403900
$ws.Cells.Item(72,13).Value = 388600
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(73,6).Value = 0
$ws.Cells.Item(73,7).Value = 0
$ws.Cells.Item(73,8).Value = 0
$ws.Cells.Item(73,9).Value = 0
$ws.Cells.Item(73,10).Value = 0
$ws.Cells.Item(73,11).Value = 0
$ws.Cells.Item(73,12).Value = 0
$ws.Cells.Item(73,13).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(74,6).Value = 0
$ws.Cells.Item(74,7).Value = 0
$ws.Cells.Item(74,8).Value = 0
$ws.Cells.Item(74,9).Value = 0
$ws.Cells.Item(74,10).Value = 0
$ws.Cells.Item(74,11).Value = 0
$ws.Cells.Item(74,12).Value = 0
$ws.Cells.Item(74,13).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(75,6).Value = 0
$ws.Cells.Item(75,7).Value = 0
$ws.Cells.Item(75,8).Value = 0
$ws.Cells.Item(75,9).Value = 0
$ws.Cells.Item(75,10).Value = 0
$ws.Cells.Item(75,11).Value = 0
$ws.Cells.Item(75,12).Value = 0
$ws.Cells.Item(75,13).Value = 0
$ws.Cells.Item(76,4).Value = 884400
$ws.Cells.Item(76,5).Value = 818900
$ws.Cells.Item(76,6).Value = 753800
$ws.Cells.Item(76,7).Value = 725900
$ws.Cells.Item(76,8).Value = 597700
$ws.Cells.Item(76,9).Value = 565200
$ws.Cells.Item(76,10).Value = 387000
$ws.Cells.Item(76,11).Value = 449100
$ws.Cells.Item(76,12).Value = 430900
$ws.Cells.Item(76,13).Value = 415700
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(77,6).Value = 0
$ws.Cells.Item(77,7).Value = 0
$ws.Cells.Item(77,8).Value = 0
$ws.Cells.Item(77,9).Value = 0
$ws.Cells.Item(77,10).Value = 0
$ws.Cells.Item(77,11).Value = 0
$ws.Cells.Item(77,12).Value = 0
$ws.Cells.Item(77,13).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(80,5).Value = 43373
$ws.Cells.Item(80,6).Value = 43281
$ws.Cells.Item(80,7).Value = 43190
$ws.Cells.Item(80,8).Value = 43100
$ws.Cells.Item(80,9).Value = 43008
$ws.Cells.Item(80,10).Value = 42916
$ws.Cells.Item(80,11).Value = 42825
$ws.Cells.Item(80,12).Value = 42735
$ws.Cells.Item(80,13).Value = 42643
$ws.Cells.Item(81,4).Value = 67400
$ws.Cells.Item(81,5).Value = 69000
$ws.Cells.Item(81,6).Value = 48700
$ws.Cells.Item(81,7).Value = 35400
$ws.Cells.Item(81,8).Value = 32600
$ws.Cells.Item(81,9).Value = 26000
$ws.Cells.Item(81,10).Value = -6500
$ws.Cells.Item(81,11).Value = 42000
$ws.Cells.Item(81,12).Value = 15200
$ws.Cells.Item(81,13).Value = 24400
$ws.Cells.Item(83,4).Value = 5700
$ws.Cells.Item(83,5).Value = 5000
$ws.Cells.Item(83,6).Value = 4800
$ws.Cells.Item(83,7).Value = 4600
$ws.Cells.Item(83,8).Value = "NA"
$ws.Cells.Item(83,9).Value = "NA"
$ws.Cells.Item(83,10).Value = "NA"
$ws.Cells.Item(83,11).Value = "NA"
$ws.Cells.Item(83,12).Value = "NA"
$ws.Cells.Item(83,13).Value = "NA"
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(84,6).Value = 0
$ws.Cells.Item(84,7).Value = 0
$ws.Cells.Item(84,8).Value = 0
$ws.Cells.Item(84,9).Value = 0
$ws.Cells.Item(84,10).Value = 0
$ws.Cells.Item(84,11).Value = 0
$ws.Cells.Item(84,12).Value = 0
$ws.Cells.Item(84,13).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(85,6).Value = 0
$ws.Cells.Item(85,7).Value = 0
$ws.Cells.Item(85,8).Value = 0
$ws.Cells.Item(85,9).Value = 0
$ws.Cells.Item(85,10).Value = 0
$ws.Cells.Item(85,11).Value = 0
$ws.Cells.Item(85,12).Value = 0
$ws.Cells.Item(85,13).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(86,6).Value = 0
$ws.Cells.Item(86,7).Value = 0
$ws.Cells.Item(86,8).Value = 0
$ws.Cells.Item(86,9).Value = 0
$ws.Cells.Item(86,10).Value = 0
$ws.Cells.Item(86,11).Value = 0
$ws.Cells.Item(86,12).Value = 0
$ws.Cells.Item(86,13).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(87,6).Value = 0
$ws.Cells.Item(87,7).Value = 0
$ws.Cells.Item(87,8).Value = 0
$ws.Cells.Item(87,9).Value = 0
$ws.Cells.Item(87,10).Value = 0
$ws.Cells.Item(87,11).Value = 0
$ws.Cells.Item(87,12).Value = 0
$ws.Cells.Item(87,13).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(88,6).Value = 0
$ws.Cells.Item(88,7).Value = 0
$ws.Cells.Item(88,8).Value = 0
$ws.Cells.Item(88,9).Value = 0
$ws.Cells.Item(88,10).Value = 0
$ws.Cells.Item(88,11).Value = 0
$ws.Cells.Item(88,12).Value = 0
$ws.Cells.Item(88,13).Value = 0
$ws.Cells.Item(89,4).Value = 440500
$ws.Cells.Item(89,5).Value = 336200
$ws.Cells.Item(89,6).Value = 62900
$ws.Cells.Item(89,7).Value = -57400
$ws.Cells.Item(89,8).Value = 214700
$ws.Cells.Item(89,9).Value = -331900
$ws.Cells.Item(89,10).Value = -144100
$ws.Cells.Item(89,11).Value = -93300
$ws.Cells.Item(89,12).Value = 510000
$ws.Cells.Item(89,13).Value = 124400
$ws.Cells.Item(91,4).Value = -37800
$ws.Cells.Item(91,5).Value = 12400
$ws.Cells.Item(91,6).Value = -7400
$ws.Cells.Item(91,7).Value = -7200
$ws.Cells.Item(91,8).Value = -6300
$ws.Cells.Item(91,9).Value = -4000
$ws.Cells.Item(91,10).Value = -4100
$ws.Cells.Item(91,11).Value = -3900
$ws.Cells.Item(91,12).Value = -15200
$ws.Cells.Item(91,13).Value = -3500
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(92,6).Value = 0
$ws.Cells.Item(92,7).Value = 0
$ws.Cells.Item(92,8).Value = 0
$ws.Cells.Item(92,9).Value = 0
$ws.Cells.Item(92,10).Value = 0
$ws.Cells.Item(92,11).Value = 0
$ws.Cells.Item(92,12).Value = 0
$ws.Cells.Item(92,13).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(93,6).Value = 0
$ws.Cells.Item(93,7).Value = 0
$ws.Cells.Item(93,8).Value = 0
$ws.Cells.Item(93,9).Value = 0
$ws.Cells.Item(93,10).Value = 0
$ws.Cells.Item(93,11).Value = 0
$ws.Cells.Item(93,12).Value = 0
$ws.Cells.Item(93,13).Value = 0
$ws.Cells.Item(94,4).Value = -37700
$ws.Cells.Item(94,5).Value = 24600
$ws.Cells.Item(94,6).Value = 0
$ws.Cells.Item(94,7).Value = -7100
$ws.Cells.Item(94,8).Value = -1200
$ws.Cells.Item(94,9).Value = -16200
$ws.Cells.Item(94,10).Value = -8900
$ws.Cells.Item(94,11).Value = -12700
$ws.Cells.Item(94,12).Value = -31700
$ws.Cells.Item(94,13).Value = -12600
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(96,6).Value = 0
$ws.Cells.Item(96,7).Value = 0
$ws.Cells.Item(96,8).Value = -2500
$ws.Cells.Item(96,9).Value = -20900
$ws.Cells.Item(96,10).Value = 0
$ws.Cells.Item(96,11).Value = 0
$ws.Cells.Item(96,12).Value = 0
$ws.Cells.Item(96,13).Value = -25900
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(97,6).Value = 0
$ws.Cells.Item(97,7).Value = 0
$ws.Cells.Item(97,8).Value = 0
$ws.Cells.Item(97,9).Value = 0
$ws.Cells.Item(97,10).Value = 0
$ws.Cells.Item(97,11).Value = 0
$ws.Cells.Item(97,12).Value = 0
$ws.Cells.Item(97,13).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(98,6).Value = 0
$ws.Cells.Item(98,7).Value = 0
$ws.Cells.Item(98,8).Value = 0
$ws.Cells.Item(98,9).Value = 0
$ws.Cells.Item(98,10).Value = 0
$ws.Cells.Item(98,11).Value = 0
$ws.Cells.Item(98,12).Value = 0
$ws.Cells.Item(98,13).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(99,6).Value = 0
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = 0
$ws.Cells.Item(99,9).Value = 0
$ws.Cells.Item(99,10).Value = 0
$ws.Cells.Item(99,11).Value = 0
$ws.Cells.Item(99,12).Value = 0
$ws.Cells.Item(99,13).Value = 0
$ws.Cells.Item(100,4).Value = -7400
$ws.Cells.Item(100,5).Value = 30800
$ws.Cells.Item(100,6).Value = 87500
$ws.Cells.Item(100,7).Value = 2900
$ws.Cells.Item(100,8).Value = 6800
$ws.Cells.Item(100,9).Value = 132800
$ws.Cells.Item(100,10).Value = 35700
$ws.Cells.Item(100,11).Value = -57400
$ws.Cells.Item(100,12).Value = 40900
$ws.Cells.Item(100,13).Value = -9700
$ws.Cells.Item(101,4).Value = -138300
$ws.Cells.Item(101,5).Value = 373800
$ws.Cells.Item(101,6).Value = 253600
$ws.Cells.Item(101,7).Value = 31500
$ws.Cells.Item(101,8).Value = "NA"
$ws.Cells.Item(101,9).Value = "NA"
$ws.Cells.Item(101,10).Value = "NA"
$ws.Cells.Item(101,11).Value = "NA"
$ws.Cells.Item(101,12).Value = "NA"
$ws.Cells.Item(101,13).Value = "NA"
$ws.Cells.Item(102,4).Value = 257200
$ws.Cells.Item(102,5).Value = 765300
$ws.Cells.Item(102,6).Value = 404000
$ws.Cells.Item(102,7).Value = -30200
$ws.Cells.Item(102,8).Value = 261700
$ws.Cells.Item(102,9).Value = -215200
$ws.Cells.Item(102,10).Value = -117300
$ws.Cells.Item(102,11).Value = -163400
$ws.Cells.Item(102,12).Value = 473300
$ws.Cells.Item(102,13).Value = 102000
